$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 135 (shifts old rows 135..222 down to 136..223)
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new weekly observation
$ws.Cells.Item(135, 1).Value = 11
$ws.Cells.Item(135, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(135, 3).Value = "Bíobío"
$ws.Cells.Item(135, 4).Value = 44873
$ws.Cells.Item(135, 5).Value = 8
$ws.Cells.Item(135, 6).Value = 100112003
$ws.Cells.Item(135, 7).Value = "Ajo"
$ws.Cells.Item(135, 8).Value = "Chino"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 350
$ws.Cells.Item(135, 11).Value = 12000
$ws.Cells.Item(135, 12).Value = 13000
$ws.Cells.Item(135, 13).Value = 12429
$ws.Cells.Item(135, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(135, 15).Value = "China"
$ws.Cells.Item(135, 16).Value = 1243
$ws.Cells.Item(135, 17).Value = 10
$ws.Cells.Item(135, 18).Value = "Hortaliza"
